$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3937473893165588
$ws.Range("B1").Value = 0.6551966071128845
$ws.Range("C1").Value = 2.214112043380737
$ws.Range("D1").Value = 4.756337642669678
$ws.Range("E1").Value = 2.091179132461548
